# Update the "last_edited_time" (column D) values exported from Notion.
# The underlying timestamps for 2024-08-03 moved from 03:54/03:55/03:56 UTC
# to 20:14/20:15/20:16 UTC, and the boundary between the "...:14" group and
# the "...:15" group shifted down from row 31 to row 15, while the boundary
# between the "...:15" group and the "...:16" group shifted down from row 85
# to row 55 (part of the multi-process run strategy change / template rework
# mentioned in the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2:D14").Value = "2024-08-03T20:14:00.000Z"
$ws.Range("D15:D54").Value = "2024-08-03T20:15:00.000Z"
$ws.Range("D55:D94").Value = "2024-08-03T20:16:00.000Z"
